# Update gh-pages to output generated at 456a3b4
# Refresh the "want to go" counts (column F) across all four sheets, and
# flip two shows whose tickets are no longer sellable (column G -> "不可售").

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 396
$ws.Range("F3").Value  = 395
$ws.Range("F4").Value  = 2643
$ws.Range("F6").Value  = 47
$ws.Range("F8").Value  = 291
$ws.Range("F10").Value = 569
$ws.Range("F11").Value = 269
$ws.Range("F12").Value = 75
$ws.Range("F13").Value = 10941
$ws.Range("F14").Value = 6268
$ws.Range("F16").Value = 14
$ws.Range("F17").Value = 402
$ws.Range("F18").Value = 241
$ws.Range("F21").Value = 868
$ws.Range("F23").Value = 213
$ws.Range("F24").Value = 882
$ws.Range("F25").Value = 3589
$ws.Range("F29").Value = 147
$ws.Range("F30").Value = 291
$ws.Range("F31").Value = 258
$ws.Range("F33").Value = 4916
$ws.Range("F35").Value = 1183
$ws.Range("F36").Value = 185
$ws.Range("F37").Value = 241
$ws.Range("F38").Value = 139

# ---- Sheet "演出" (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value  = "不可售"
$ws.Range("F13").Value = 3632
$ws.Range("F21").Value = 78
$ws.Range("F22").Value = 12
$ws.Range("F24").Value = 7

# ---- Sheet "本地生活" (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8922
$ws.Range("F3").Value = 461
$ws.Range("F4").Value = 1727

# ---- Sheet "全部类型" (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 461
$ws.Range("F3").Value  = 1727
$ws.Range("G4").Value  = "不可售"
$ws.Range("F6").Value  = 395
$ws.Range("F7").Value  = 2643
$ws.Range("F13").Value = 47
$ws.Range("F14").Value = 291
$ws.Range("F17").Value = 269
$ws.Range("F18").Value = 10941
$ws.Range("F19").Value = 3632
$ws.Range("F22").Value = 14
$ws.Range("F23").Value = 402
$ws.Range("F24").Value = 241
$ws.Range("F28").Value = 213
$ws.Range("F29").Value = 882
$ws.Range("F30").Value = 3589
$ws.Range("F33").Value = 147
$ws.Range("F34").Value = 291
$ws.Range("F35").Value = 258
$ws.Range("F40").Value = 4916
$ws.Range("F42").Value = 1183
$ws.Range("F44").Value = 185
$ws.Range("F45").Value = 139
$ws.Range("F47").Value = 78
$ws.Range("F49").Value = 7
